# Update resistor/component placement values on the DuDad-bottom-pos sheet
# and remove the stray empty "Sheet1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DuDad-bottom-pos")

# J1 rotation: 90 -> 0
$ws.Range("D2").Value = 0

# J3 Mid X / Mid Y shift
$ws.Range("B3").Value = 55.625
$ws.Range("C3").Value = -123.325

# Remove the empty, unused "Sheet1" worksheet
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()
